# RPAR_holdings.xlsx update — "Add files via upload"
# Refresh the model-holdings snapshot date and the Weight/Percent Change
# figures that go with it, then restore the original sheet protection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet ships protected; unprotect it so the cells can be written, then
# re-protect it once the data refresh is done.
$ws.Unprotect()

# Disclosure footer: bump the "as of" date.
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-29 for illustrative purposes only and are subject to change."

# Refreshed Weight (col D) / Percent Change (col E) figures for rows 2-15.
$ws.Range("D2").Value  = 0.0557029455019912
$ws.Range("E2").Value  = -0.00469801908267542

$ws.Range("D3").Value  = 0.02339391515397593
$ws.Range("E3").Value  = -0.004843592330978641

$ws.Range("D4").Value  = 0.03150439027943537
$ws.Range("E4").Value  = -0.005587668593448836

$ws.Range("D5").Value  = 0.03230633247687165
$ws.Range("E5").Value  = -0.005950019833399489

$ws.Range("D6").Value  = 0.03386367859451705
$ws.Range("E6").Value  = -0.006214149139579517

$ws.Range("D7").Value  = 0.01891578628948629
$ws.Range("E7").Value  = -0.004867090977162092

$ws.Range("D8").Value  = 0.004645734109286038
$ws.Range("E8").Value  = -0.02395470383275244

$ws.Range("D9").Value  = 0.006650252368983789
$ws.Range("E9").Value  = -0.004868154158214955

$ws.Range("D10").Value = 0.06960507550498246
$ws.Range("E10").Value = -0.01162790697674421

$ws.Range("D11").Value = 0.06976694777359868
$ws.Range("E11").Value = -0.01218097447795807

$ws.Range("D12").Value = 0.1474764281939675
$ws.Range("E12").Value = -0.008488218937509173

$ws.Range("D13").Value = 0.3906908607254373
$ws.Range("E13").Value = -0.003433703116745956

$ws.Range("D14").Value = 0.1154776530274667
$ws.Range("E14").Value = -0.008585796639273946

$ws.Range("D15").Value = 0.9999999999999999
$ws.Range("E15").Value = -0.006433410975816733

# Restore sheet protection (workbook shipped protected, no password changes
# described by this edit).
$ws.Protect()
